$wb = $excel.ActiveWorkbook

# --- Fix the "maxDetxerity" -> "maxDexterity" typo in the shared strings ---
# Every sheet (human, dwarf, elf) has the same header layout, with column H
# holding this label on row 1. Re-writing the cell value on every sheet drops
# the old (now-unused) shared string and appends the corrected spelling as a
# brand new shared-string entry, matching the diff (old entry removed from
# its old slot, new entry appended at the end of the table).
foreach ($ws in $wb.Worksheets) {
    $ws.Range("H1").Value = "maxDexterity"
}

# --- human sheet: update the view (zoom + selection) ---
$wsHuman = $wb.Worksheets.Item(1)
$wsHuman.Activate() | Out-Null
$wsHuman.Range("J10").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85

# --- dwarf sheet: update the selection ---
$wsDwarf = $wb.Worksheets.Item(2)
$wsDwarf.Activate() | Out-Null
$wsDwarf.Range("I5").Select() | Out-Null

# --- elf sheet: update the selection ---
$wsElf = $wb.Worksheets.Item(3)
$wsElf.Activate() | Out-Null
$wsElf.Range("H3").Select() | Out-Null

# Leave the human sheet active/selected, matching tabSelected="1" on sheet1.
$wsHuman.Activate() | Out-Null
